# Add season-record columns (Wins, Losses, Ties) as AD, AE, AF.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): set the labels, then copy the existing header
# formatting (bold, centered, bordered) from AB1 onto the new header cells
# so they pick up the same style used by the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AB1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2-51: every player row gets the same season-record values.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 72  # AD - Wins
    $ws.Cells.Item($row, 31).Value = 90  # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF - Ties
}
